# Add team record (Wins/Losses/Ties) columns to the BAL_1994 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new column headers in AD1:AF1
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style (bold, centered, bordered) from an existing header cell
# onto the new header cells so they match the rest of row 1.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Data rows 2-34: every row gets the same team record values.
for ($r = 2; $r -le 34; $r++) {
    $ws.Cells.Item($r, 30).Value = 63  # AD = Wins
    $ws.Cells.Item($r, 31).Value = 49  # AE = Losses
    $ws.Cells.Item($r, 32).Value = 0   # AF = Ties
}

Write-Host "Added Wins/Losses/Ties columns (AD:AF) for rows 1-34"
